$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.19"
$ws.Range("E2").Value = "'1.04%"
$ws.Range("D3").Value = "'44.13"
$ws.Range("E3").Value = "'-1.06%"
$ws.Range("D4").Value = "'5.495"
$ws.Range("E4").Value = "'-2.11%"
$ws.Range("D5").Value = "'0.08011"
$ws.Range("E5").Value = "'-0.54%"
$ws.Range("D6").Value = "'1.977"
$ws.Range("E6").Value = "'4.03%"
$ws.Range("B7").Value = "'MXToken"
$ws.Range("C7").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9525"
$ws.Range("E7").Value = "'0.82%"
$ws.Range("B8").Value = "'BTSEToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.551"
$ws.Range("E8").Value = "'-3.97%"
$ws.Range("D9").Value = "'0.1137"
$ws.Range("E9").Value = "'-2.40%"
$ws.Range("D10").Value = "'0.1889"
$ws.Range("E10").Value = "'2.01%"
$ws.Range("D11").Value = "'10.70"
$ws.Range("E11").Value = "'27.45%"
$ws.Range("D12").Value = "'0.09958"
$ws.Range("E12").Value = "'0.78%"
$ws.Range("D13").Value = "'0.04819"
$ws.Range("E13").Value = "'13.30%"
$ws.Range("D14").Value = "'0.1064"
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("D15").Value = "'0.001270"
$ws.Range("E15").Value = "'-0.81%"
$ws.Range("D16").Value = "'0.04070"
$ws.Range("E16").Value = "'-3.52%"
$ws.Range("D17").Value = "'0.005922"
$ws.Range("E17").Value = "'0.62%"
$ws.Range("B18").Value = "'HotbitToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D18").Value = "'0.004348"
$ws.Range("E18").Value = "'-2.53%"
$ws.Range("B19").Value = "'LEO"
$ws.Range("C19").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.372"
$ws.Range("E19").Value = "'-6.37%"
$ws.Range("B20").Value = "'GateToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D20").Value = "'4.385"
$ws.Range("E20").Value = "'1.95%"
$ws.Range("B21").Value = "'BitpandaEcosystemToken"
$ws.Range("C21").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "'0.3484"
$ws.Range("E21").Value = "'-0.39%"
$ws.Range("B22").Value = "'ProBitToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1415"
$ws.Range("E22").Value = "'3.21%"
$ws.Range("B23").Value = "'ZBToken"
$ws.Range("C23").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.2584"
$ws.Range("E23").Value = "'-1.06%"
$ws.Range("B24").Value = "'BitKan"
$ws.Range("C24").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001270"
$ws.Range("E24").Value = "'1.98%"
$ws.Range("D25").Value = "'0.0001198"
$ws.Range("E25").Value = "'1.35%"
$ws.Range("D26").Value = "'0.0003738"
$ws.Range("E26").Value = "'-6.46%"
$ws.Range("D38").Value = "'0.02603"
$ws.Range("E38").Value = "'-1.16%"
$ws.Range("D39").Value = "'0.05748"
$ws.Range("E39").Value = "'4.77%"
$ws.Range("D40").Value = "'0.007542"
$ws.Range("E40").Value = "'-1.33%"
$ws.Range("D41").Value = "'0.1403"
$ws.Range("E41").Value = "'0.52%"
$ws.Range("D42").Value = "'0.007342"
$ws.Range("E42").Value = "'2.87%"
$ws.Range("D43").Value = "'0.002012"
$ws.Range("E43").Value = "'0.41%"
$ws.Range("D44").Value = "'0.008794"
$ws.Range("E44").Value = "'-0.70%"
$ws.Range("D45").Value = "'0.00007117"
$ws.Range("E45").Value = "'-0.06%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.34%"
$ws.Range("D47").Value = "'0.0005794"
$ws.Range("E47").Value = "'-0.31%"
$ws.Range("D48").Value = "'0.003524"
$ws.Range("E48").Value = "'54.98%"
$ws.Range("E49").Value = "'-3.60%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.34%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.34%"
